$d = $word.ActiveDocument

# Locate the CGPA text in the document.
$old = "CGPA: 7.33/10"
$before = "CGPA: 7.3"
$digit = "4"

$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target CGPA text"
}

$s = $rng.Start

# The single character between "CGPA: 7.3" and "/10" is the digit that
# changed (was "3", becomes "4"). Replace just that character.
$midStart = $s + $before.Length
$midEnd = $midStart + 1
$mid = $d.Range($midStart, $midEnd)
$mid.Text = $digit

# Force the run containing the replaced digit to split away from its
# neighbours (which keep identical formatting) by toggling a character
# property on just that sub-range and reverting it. This mirrors how a
# targeted, in-place edit (rather than a full-run rewrite) leaves the
# paragraph with three runs: "CGPA: 7.3", "4", "/10".
$mid2 = $d.Range($midStart, $midEnd)
$mid2.Font.Bold = 1
$mid2.Font.Bold = 0
